$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'246.77"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Formula = "'21.55"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").Formula = "'5.284"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Formula = "'0.05583"
$ws.Range("D5").ClearFormats()
$ws.Range("D7").Formula = "'6.389"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Formula = "'0.8168"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Formula = "'0.9710"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Formula = "'0.1403"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Formula = "'0.07492"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Formula = "'0.03145"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Formula = "'0.03049"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Formula = "'0.09290"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Formula = "'3.582"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Formula = "'0.001587"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Formula = "'0.04717"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Formula = "'0.0005774"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Formula = "'0.006372"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Formula = "'0.005050"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Formula = "'0.001031"
$ws.Range("D21").ClearFormats()
$ws.Range("D23").Formula = "'3.779"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").Formula = "'0.3257"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Formula = "'0.1286"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").Formula = "'0.0003096"
$ws.Range("D28").ClearFormats()
$ws.Range("D40").Formula = "'0.03928"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Formula = "'0.007026"
$ws.Range("D41").ClearFormats()
$ws.Range("D43").Formula = "'0.003398"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Formula = "'0.007801"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Formula = "'0.00005802"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Formula = "'0.00000000750"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").Formula = "'0.0005494"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Formula = "'0.6792"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Formula = "'0.1451"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Formula = "'0.00002099"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").Formula = "'0.01009"
$ws.Range("D51").ClearFormats()
